$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 7 (2021年) - copy row 6's formatting first so the label cell
# keeps the same style (bold, centered, bordered) as the other year cells.
$ws.Range("A6:H6").Copy()
$ws.Range("A7:H7").PasteSpecial(-4122)
$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 103.9
$ws.Range("C7").Value = 101.2
$ws.Range("D7").Value = 117.3
$ws.Range("E7").Value = 99.7
$ws.Range("F7").Value = 105.2
$ws.Range("G7").Value = 101.4
$ws.Range("H7").Value = 100.7

# Add row 8 (2022年) - only B8 has data so far; C8:H8 stay blank
# (still present as touched/empty cells, matching the source data file).
$ws.Range("A6:H6").Copy()
$ws.Range("A8:H8").PasteSpecial(-4122)
$ws.Range("A8").Value = "2022年"
$ws.Range("B8").Value = 105
